$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.449521780014038
$ws.Range("B1").Value = 2.794248342514038
$ws.Range("C1").Value = 1.673259258270264
$ws.Range("D1").Value = 1.358235716819763
$ws.Range("E1").Value = 1.266934752464294
